$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "304.64"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.43%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.12"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.64%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.030"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.36%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07822"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.08%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.171"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.15%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.908"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.12%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9178"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.44%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09739"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.97%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1859"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.42%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08666"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.90%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03478"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.74%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09911"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.11%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001445"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.28%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005675"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.26%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.460"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.34%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.093"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.52%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.371"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "15.56%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3423"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.15%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1345"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.34%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.762"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.50%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2207"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.95%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04596"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.00%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "14.80%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001229"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.09%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001399"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "7.77%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004750"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.12%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01825"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.63%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04727"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.00%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007632"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.95%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1394"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.51%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007742"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.30%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002238"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.11%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01111"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "10.29%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006377"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.65%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.13%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0005800"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.01%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.63"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "183.56%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.002000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-25.89%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00002100"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.13%"
